# Auto-generated edit script applying market-data refresh to Phantom_Profits workbook
# Updates columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*) per scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 293.42856
$ws.Range("I8").Value = 293.42856
$ws.Range("K8").Value = 880.28568
$ws.Range("M8").Value = -741.28568

$ws.Range("H11").Value = 21.428572
$ws.Range("I11").Value = 21.428572
$ws.Range("K11").Value = 21.428572
$ws.Range("M11").Value = 118.571428

$ws.Range("H31").Value = 45
$ws.Range("I31").Value = 45
$ws.Range("K31").Value = 135
$ws.Range("M31").Value = 95

$ws.Range("H33").Value = 389.83334
$ws.Range("I33").Value = 407.9
$ws.Range("K33").Value = 407.9
$ws.Range("M33").Value = -178.9

$ws.Range("H38").Value = 272.22223
$ws.Range("I38").Value = 272.22223
$ws.Range("K38").Value = 816.66669
$ws.Range("M38").Value = -444.66669

$ws.Range("H39").Value = 215.22223
$ws.Range("I39").Value = 179.625
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 538.875
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -242.875
$ws.Range("N39").Value = -2092

$ws.Range("H131").Value = 2195.4
$ws.Range("I131").Value = 2195.4
$ws.Range("K131").Value = 6586.200000000001
$ws.Range("M131").Value = -1546.200000000001

$ws.Range("H132").Value = 3866
$ws.Range("I132").Value = 3866
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11598
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9068
$ws.Range("N132").ClearContents()

$ws.Range("H138").Value = 2608.0667
$ws.Range("I138").Value = 2608.0667
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 7824.2001
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -2684.2001
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1072.5
$ws.Range("I2").Value = 1072.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1072.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -959.5
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 14937.619
$ws.Range("I32").Value = 14484.3
$ws.Range("K32").Value = 14484.3
$ws.Range("M32").Value = -14197.3

$ws.Range("H45").Value = 3999.5
$ws.Range("I45").Value = 3999.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3999.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3622.5
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 10541.158
$ws.Range("I61").Value = 10877.857
$ws.Range("J61").Value = 9598.4
$ws.Range("K61").Value = 10877.857
$ws.Range("L61").Value = 9598.4
$ws.Range("M61").Value = -10665.857
$ws.Range("N61").Value = -10022.4

$ws.Range("H116").Value = 1072.5
$ws.Range("I116").Value = 1072.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1072.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1221.5
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 10541.158
$ws.Range("I136").Value = 10877.857
$ws.Range("J136").Value = 9598.4
$ws.Range("K136").Value = 32633.571
$ws.Range("L136").Value = 28795.2
$ws.Range("M136").Value = -30083.571
$ws.Range("N136").Value = -33895.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1072.5
$ws.Range("I3").Value = 1072.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1072.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -958.5
$ws.Range("N3").ClearContents()

$ws.Range("H107").Value = 2750
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -7840

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 5857.143
$ws.Range("I134").Value = 5704.5
$ws.Range("K134").Value = 17113.5
$ws.Range("M134").Value = -14578.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3424.75
$ws.Range("I16").Value = 1233.3334
$ws.Range("K16").Value = 1233.3334
$ws.Range("M16").Value = -946.3334

$ws.Range("H31").Value = 5115
$ws.Range("I31").Value = 3109.8572
$ws.Range("J31").Value = 9793.666999999999
$ws.Range("K31").Value = 3109.8572
$ws.Range("L31").Value = 9793.666999999999
$ws.Range("M31").Value = -2814.8572
$ws.Range("N31").Value = -10383.667

$ws.Range("H34").Value = 5115
$ws.Range("I34").Value = 3109.8572
$ws.Range("J34").Value = 9793.666999999999
$ws.Range("K34").Value = 3109.8572
$ws.Range("L34").Value = 9793.666999999999
$ws.Range("M34").Value = -2907.8572
$ws.Range("N34").Value = -10197.667

$ws.Range("H58").Value = 2959.6875
$ws.Range("I58").Value = 1965.8462
$ws.Range("K58").Value = 1965.8462
$ws.Range("M58").Value = -1762.8462

$ws.Range("H113").Value = 3424.75
$ws.Range("I113").Value = 1233.3334
$ws.Range("K113").Value = 1233.3334
$ws.Range("M113").Value = 936.6666

$ws.Range("H136").Value = 2959.6875
$ws.Range("I136").Value = 1965.8462
$ws.Range("K136").Value = 5897.5386
$ws.Range("M136").Value = -3347.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2637.1538
$ws.Range("I5").Value = 2637.1538
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 7911.4614
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -7799.4614
$ws.Range("N5").ClearContents()

$ws.Range("H39").Value = 6169.5264
$ws.Range("J39").Value = 6818.9414
$ws.Range("L39").Value = 20456.8242
$ws.Range("N39").Value = -21044.8242

$ws.Range("H135").Value = 2637.1538
$ws.Range("I135").Value = 2637.1538
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 23734.3842
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -21199.3842
$ws.Range("N135").ClearContents()

$ws.Range("H139").Value = 7077.8
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10692.23
$ws.Range("J43").Value = 22999.8
$ws.Range("L43").Value = 22999.8
$ws.Range("N43").Value = -23301.8

$ws.Range("H122").Value = 2394.3635
$ws.Range("I122").Value = 2394.3635
$ws.Range("K122").Value = 7183.0905
$ws.Range("M122").Value = -4733.0905

$ws.Range("H126").Value = 2608
$ws.Range("J126").Value = 1200
$ws.Range("L126").Value = 3600
$ws.Range("N126").Value = -8540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5365.778
$ws.Range("I7").Value = 5424
$ws.Range("K7").Value = 5424
$ws.Range("M7").Value = -5312

$ws.Range("H22").Value = 749.5
$ws.Range("I22").Value = 666
$ws.Range("K22").Value = 666
$ws.Range("M22").Value = -371

$ws.Range("H27").Value = 749.5
$ws.Range("I27").Value = 666
$ws.Range("K27").Value = 666
$ws.Range("M27").Value = -559

$ws.Range("H93").Value = 1383.421
$ws.Range("I93").Value = 1498.2307
$ws.Range("J93").Value = 1134.6666
$ws.Range("K93").Value = 1498.2307
$ws.Range("L93").Value = 1134.6666
$ws.Range("M93").Value = -250.2307000000001
$ws.Range("N93").Value = -3630.6666

$ws.Range("H122").Value = 2384.4285
$ws.Range("I122").Value = 2233
$ws.Range("K122").Value = 6699
$ws.Range("M122").Value = -4249

$ws.Range("H126").Value = 5365.778
$ws.Range("I126").Value = 5424
$ws.Range("K126").Value = 16272
$ws.Range("M126").Value = -13802

$ws.Range("H136").Value = 7619.1763
$ws.Range("I136").Value = 6109
$ws.Range("K136").Value = 18327
$ws.Range("M136").Value = -15777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H81").Value = 5054.778
$ws.Range("I81").Value = 5286.625
$ws.Range("K81").Value = 10573.25
$ws.Range("M81").Value = -9512.25

$ws.Range("H84").Value = 5054.778
$ws.Range("I84").Value = 5286.625
$ws.Range("K84").Value = 52866.25
$ws.Range("M84").Value = -47562.25

$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492

$ws.Range("H107").Value = 843.2
$ws.Range("J107").Value = 577.6667
$ws.Range("L107").Value = 1733.0001
$ws.Range("N107").Value = -5573.0001

$ws.Range("H113").Value = 691.2
$ws.Range("I113").Value = 569.0909
$ws.Range("J113").Value = 840.44446
$ws.Range("K113").Value = 1707.2727
$ws.Range("L113").Value = 2521.33338
$ws.Range("M113").Value = 462.7273
$ws.Range("N113").Value = -6861.33338

$ws.Range("H126").Value = 898.4545000000001
$ws.Range("I126").Value = 738.3
$ws.Range("K126").Value = 2214.9
$ws.Range("M126").Value = 255.1000000000004

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
